$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (shifts existing rows 26..45 down to 27..46,
# row 45's old data ends up as row 46, matching target dimension A1:R46).
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the weekly record.
$ws.Cells.Item(26, 1).Value = 9
$ws.Cells.Item(26, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(26, 3).Value = "Metropolitana"
$ws.Cells.Item(26, 4).Value = 45001
$ws.Cells.Item(26, 5).Value = 13
$ws.Cells.Item(26, 6).Value = 100112010
$ws.Cells.Item(26, 7).Value = "Achicoria"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 90
$ws.Cells.Item(26, 11).Value = 6000
$ws.Cells.Item(26, 12).Value = 6000
$ws.Cells.Item(26, 13).Value = 6000
$ws.Cells.Item(26, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(26, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(26, 16).Value = 375
$ws.Cells.Item(26, 17).Value = 16
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Row 36's Origen differs between before/after states: it used to read
# "Región Metropolitana" in the old row 36, which is now (post-shift) row 37.
# The target keeps "Provincia de Quillota" for the new row 36 (old row 35)
# and "Región Metropolitana" for row 37 (old row 36) -- both already carried
# over correctly by the row insert/shift, so no further action is needed here.
